$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: refresh the scrape timestamp (column A) on the existing rows 2-13; nothing else in those rows changed
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("A$r").Value = '2025-11-12 01:49:51'
}

# Step 2: drop all existing hyperlinks; they will be rebuilt below once the cell text has its final values/positions
$ws.Hyperlinks.Delete()

# Step 3: write rows 14-23.
#   Row 14 is a brand-new scraped listing.
#   Rows 15-23 are the former rows 14-22, each shifted down by one row.
# Row 14
$ws.Range("A14").Value = '2025-11-12 01:49:51'
$ws.Range("B14").Value = '【急募】楽天市場在庫連動システム(同一店舗内)のエラー修正依頼'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5432212'
$ws.Range("G14").Value = 25

# Row 15
$ws.Range("A15").Value = '2025-11-12 01:49:51'
$ws.Range("B15").Value = '【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5431911'
$ws.Range("G15").Value = 25

# Row 16
$ws.Range("A16").Value = '2025-11-12 01:49:51'
$ws.Range("B16").Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Range("G16").Value = 25

# Row 17
$ws.Range("A17").Value = '2025-11-12 01:49:51'
$ws.Range("B17").Value = '【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5417644'
$ws.Range("G17").Value = 25

# Row 18
$ws.Range("A18").Value = '2025-11-12 01:49:51'
$ws.Range("B18").Value = '【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5431322'
$ws.Range("G18").Value = 25

# Row 19
$ws.Range("A19").Value = '2025-11-12 01:49:51'
$ws.Range("B19").Value = '【音楽制作】サイケデリックトランスのトラックを作成してくれる方募集'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5432042'
$ws.Range("G19").Value = 18

# Row 20
$ws.Range("A20").Value = '2025-11-12 01:49:51'
$ws.Range("B20").Value = '初回 Hubspot構築者募集'
$ws.Range("C20").Value = 'システム開発'
$ws.Range("D20").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E20").Value = '期限情報なし'
$ws.Range("F20").Value = 'https://www.lancers.jp/work/detail/5431947'
$ws.Range("G20").Value = 18

# Row 21
$ws.Range("A21").Value = '2025-11-12 01:49:51'
$ws.Range("B21").Value = 'AWS環境からAWS環境ヘの新規構築'
$ws.Range("C21").Value = 'システム開発'
$ws.Range("D21").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E21").Value = '期限情報なし'
$ws.Range("F21").Value = 'https://www.lancers.jp/work/detail/5431069'
$ws.Range("G21").Value = 18

# Row 22
$ws.Range("A22").Value = '2025-11-12 01:49:51'
$ws.Range("B22").Value = '【Stable Diffusion】参考動画に沿って約100プロンプト構築'
$ws.Range("C22").Value = 'システム開発'
$ws.Range("D22").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E22").Value = '期限情報なし'
$ws.Range("F22").Value = 'https://www.lancers.jp/work/detail/5432055'
$ws.Range("G22").Value = 10

# Row 23
$ws.Range("A23").Value = '2025-11-12 01:49:51'
$ws.Range("B23").Value = 'EAの作成'
$ws.Range("C23").Value = 'システム開発'
$ws.Range("D23").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E23").Value = '期限情報なし'
$ws.Range("F23").Value = 'https://www.lancers.jp/work/detail/5431276'
$ws.Range("G23").Value = 10

# Step 4: (re)create the hyperlinks for F2:F23, pointing at the same literal URL each cell displays
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5431738')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5431917')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5431299')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5431740')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5431673')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5411585')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5431547')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5431786')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5432161')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5418064')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5431852')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5431508')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5432212')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5431911')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5341051')
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5417644')
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5431322')
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5432042')
$ws.Hyperlinks.Add($ws.Range("F20"), 'https://www.lancers.jp/work/detail/5431947')
$ws.Hyperlinks.Add($ws.Range("F21"), 'https://www.lancers.jp/work/detail/5431069')
$ws.Hyperlinks.Add($ws.Range("F22"), 'https://www.lancers.jp/work/detail/5432055')
$ws.Hyperlinks.Add($ws.Range("F23"), 'https://www.lancers.jp/work/detail/5431276')

# Step 5: Hyperlinks.Add() stamps a freshly duplicated xf on the cell; reapply the shared "Hyperlink" style
#         so every F-cell keeps using the single original style record (s="1") like the rest of the column.
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("F$r").Style = "Hyperlink"
}
